# Regenerate the handback-status report: the previously-seen source file
# "d6a5534b-124a-4535-8f3c-cfb62ed6460e.md" was re-handed-back under a new
# name ("10bc2547-...") and a brand-new source file
# ("e64890e6-a884-4969-9091-289243af2a76.md") was handed back as well.
# Add a row for the new file to every table (Overview, zh-cn, de-de) and
# refresh the timestamps / generated xliff names for the renamed file.

$wb = $excel.ActiveWorkbook

$renamedName = "10bc2547-d7b6-4e32-9693-5ebe9ed9adb8"
$newName = "e64890e6-a884-4969-9091-289243af2a76"

$renamedHash = "73d411c1a207f0936430d519db4b1b97ddc13935"
$newHash = "0aba362954cf6245801076894750dc3228f6aed3"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview": rename the existing row's file, then append a row
# for the newly handed-back file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = ($renamedName + ".md")
$wsOverview.Range("B2").Value = ("e2e\" + $renamedName + ".md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$renamedName.md", "", "", ("e2e\" + $renamedName + ".md")) | Out-Null
$wsOverview.Range("G2").Value = "2017-01-03 05:29:22"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = ($newName + ".md")
$wsOverview.Range("B3").Value = ("e2e\" + $newName + ".md")
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2017-01-03 05:28:52"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$newName.md", "", "", ("e2e\" + $newName + ".md")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("J2").Hyperlinks.Delete()

$wsZhCn.Range("A2").Value = ($renamedName + ".md")
$wsZhCn.Range("G2").Value = ("$renamedName.$renamedHash.zh-cn.xlf")
$wsZhCn.Range("H2").Value = "2017-01-03 05:29:12"
$wsZhCn.Range("H2").NumberFormat = $dateFmt
$wsZhCn.Range("J2").Value = ($renamedName + ".md")
$wsZhCn.Range("K2").Value = ("$renamedName.$renamedHash.zh-cn.xlf")
$wsZhCn.Range("L2").Value = "2017-01-03 05:29:38"
$wsZhCn.Range("L2").NumberFormat = $dateFmt

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$renamedName.md", "", "", ($renamedName + ".md")) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test1-zhcn/blob/3004e0498337170150356bebb5d11e156956c4fd/e2e/$renamedName.md", "", "", ($renamedName + ".md")) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = ($newName + ".md")
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = ("$newName.$newHash.zh-cn.xlf")
$wsZhCn.Range("H3").Value = "2017-01-03 05:28:41"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("J3").Value = ($newName + ".md")
$wsZhCn.Range("K3").Value = ("$newName.$newHash.zh-cn.xlf")
$wsZhCn.Range("L3").Value = "2017-01-03 05:29:38"
$wsZhCn.Range("L3").NumberFormat = $dateFmt
$wsZhCn.Range("O3").Value = "True"
$wsZhCn.Range("Q3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$newName.md", "", "", ($newName + ".md")) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test1-zhcn/blob/3004e0498337170150356bebb5d11e156956c4fd/e2e/$newName.md", "", "", ($newName + ".md")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("J2").Hyperlinks.Delete()

$wsDeDe.Range("A2").Value = ($renamedName + ".md")
$wsDeDe.Range("G2").Value = ("$renamedName.$renamedHash.de-de.xlf")
$wsDeDe.Range("H2").Value = "2017-01-03 05:29:22"
$wsDeDe.Range("H2").NumberFormat = $dateFmt
$wsDeDe.Range("J2").Value = ($renamedName + ".md")
$wsDeDe.Range("K2").Value = ("$renamedName.$renamedHash.de-de.xlf")
$wsDeDe.Range("L2").Value = "2017-01-03 05:29:50"
$wsDeDe.Range("L2").NumberFormat = $dateFmt

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$renamedName.md", "", "", ($renamedName + ".md")) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test1-dede/blob/1cf8727aa35765a6381956c755234b9cfb6bb629/e2e/$renamedName.md", "", "", ($renamedName + ".md")) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = ($newName + ".md")
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = ("$newName.$newHash.de-de.xlf")
$wsDeDe.Range("H3").Value = "2017-01-03 05:28:52"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("J3").Value = ($newName + ".md")
$wsDeDe.Range("K3").Value = ("$newName.$newHash.de-de.xlf")
$wsDeDe.Range("L3").Value = "2017-01-03 05:29:50"
$wsDeDe.Range("L3").NumberFormat = $dateFmt
$wsDeDe.Range("O3").Value = "True"
$wsDeDe.Range("Q3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/9bbfdef86382df90f0c48f5f9b3627ba4eafbd27/e2e/$newName.md", "", "", ($newName + ".md")) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test1-dede/blob/1cf8727aa35765a6381956c755234b9cfb6bb629/e2e/$newName.md", "", "", ($newName + ".md")) | Out-Null

Write-Host "Applied handback status update."
